# Daily attendance processing - 2025-10-15 11:20:27
# Reorders the "Recorded By" (column G) list so that the literal entry
# "System" is moved to the front of the comma-separated list, keeping
# the relative order of the remaining entries unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -like "*,*" -and $value -like "*System*") {
        $parts = $value -split ",\s*"
        $idx = -1
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($parts[$i].Equals("System")) {
                $idx = $i
                break
            }
        }
        if ($idx -ge 0) {
            $newParts = New-Object System.Collections.ArrayList
            [void]$newParts.Add($parts[$idx])
            for ($i = 0; $i -lt $parts.Count; $i++) {
                if ($i -ne $idx) {
                    [void]$newParts.Add($parts[$i])
                }
            }
            $cell.Value = [string]::Join(", ", $newParts)
        }
    }
}
